$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.032830953598022
$ws.Range("B1").Value = 3.426476955413818
$ws.Range("C1").Value = 3.533007621765137
$ws.Range("D1").Value = 2.039578914642334
$ws.Range("E1").Value = 1.17303478717804
